# Update the "Login" sheet (panel configuration settings):
# add a new row for an additional item code (TB7SX1CC) below the
# existing ItemCode/Quantity row, matching the style of the row above,
# and move the selection onto F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# New panel setting row: F3 holds the additional item code, formatted
# like the existing F2 cell (text number format).
$ws.Range("F3").Value = "TB7SX1CC"
$ws.Range("F3").NumberFormat = $ws.Range("F2").NumberFormat

# Update the active selection to reflect the new panel focus cell.
$ws.Range("F2").Select()
